# Auto-generated script applying numeric cell updates across multiple sheets
# per the commit diff (scheduled market-data refresh for Gilgamesh_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1626.7778
$ws.Range("J17").Value = 1767.75
$ws.Range("L17").Value = 5303.25
$ws.Range("N17").Value = -5639.25
$ws.Range("H80").Value = 577.5
$ws.Range("J80").Value = 752.5
$ws.Range("L80").Value = 2257.5
$ws.Range("N80").Value = -4253.5
$ws.Range("H83").Value = 577.5
$ws.Range("J83").Value = 752.5
$ws.Range("L83").Value = 6772.5
$ws.Range("N83").Value = -16756.5
$ws.Range("H100").Value = 1746.5714
$ws.Range("I100").Value = 1875.3334
$ws.Range("J100").Value = 1650
$ws.Range("K100").Value = 1875.3334
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -1334.3334
$ws.Range("N100").Value = -2732
$ws.Range("H125").Value = 935.2857
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H137").Value = 1194288.6
$ws.Range("I137").Value = 1788483.6
$ws.Range("J137").Value = 5898.5
$ws.Range("K137").Value = 5365450.800000001
$ws.Range("L137").Value = 17695.5
$ws.Range("M137").Value = -5362900.800000001
$ws.Range("N137").Value = -22795.5
$ws.Range("H138").Value = 2353.96
$ws.Range("J138").Value = 3058.25
$ws.Range("L138").Value = 9174.75
$ws.Range("N138").Value = -19454.75
$ws.Range("H141").Value = 3134.9333
$ws.Range("I141").Value = 2917.4167
$ws.Range("K141").Value = 8752.250100000001
$ws.Range("M141").Value = -3572.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 299.5
$ws.Range("I4").Value = 299.5
$ws.Range("K4").Value = 299.5
$ws.Range("M4").Value = -183.5
$ws.Range("H22").Value = 1999
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H32").Value = 1629147.1
$ws.Range("I32").Value = 757030.4399999999
$ws.Range("J32").Value = 15873721
$ws.Range("K32").Value = 757030.4399999999
$ws.Range("L32").Value = 15873721
$ws.Range("M32").Value = -756743.4399999999
$ws.Range("N32").Value = -15874295
$ws.Range("H61").Value = 3585.9048
$ws.Range("I61").Value = 1553.7333
$ws.Range("K61").Value = 1553.7333
$ws.Range("M61").Value = -1341.7333
$ws.Range("H74").Value = 148184.23
$ws.Range("I74").Value = 223777.12
$ws.Range("J74").Value = 2813.3076
$ws.Range("K74").Value = 223777.12
$ws.Range("L74").Value = 2813.3076
$ws.Range("M74").Value = -222903.12
$ws.Range("N74").Value = -4561.3076
$ws.Range("H77").Value = 148184.23
$ws.Range("I77").Value = 223777.12
$ws.Range("J77").Value = 2813.3076
$ws.Range("K77").Value = 1118885.6
$ws.Range("L77").Value = 14066.538
$ws.Range("M77").Value = -1114517.6
$ws.Range("N77").Value = -22802.538
$ws.Range("H80").Value = 46733.332
$ws.Range("I80").Value = 20100
$ws.Range("K80").Value = 20100
$ws.Range("M80").Value = -19102
$ws.Range("H83").Value = 46733.332
$ws.Range("I83").Value = 20100
$ws.Range("K83").Value = 60300
$ws.Range("M83").Value = -55308
$ws.Range("H110").Value = 3000
$ws.Range("H122").Value = 2832.6875
$ws.Range("J122").Value = 3875
$ws.Range("L122").Value = 11625
$ws.Range("N122").Value = -16525
$ws.Range("H136").Value = 3585.9048
$ws.Range("I136").Value = 1553.7333
$ws.Range("K136").Value = 4661.199900000001
$ws.Range("M136").Value = -2111.199900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 856.8333
$ws.Range("I22").Value = 998.5
$ws.Range("K22").Value = 998.5
$ws.Range("M22").Value = -825.5
$ws.Range("H80").Value = 895.6667
$ws.Range("I80").Value = 700
$ws.Range("K80").Value = 700
$ws.Range("M80").Value = 298
$ws.Range("H83").Value = 895.6667
$ws.Range("I83").Value = 700
$ws.Range("K83").Value = 3500
$ws.Range("M83").Value = 1492
$ws.Range("H86").Value = 5249.25
$ws.Range("I86").Value = 3999
$ws.Range("J86").Value = 6499.5
$ws.Range("K86").Value = 3999
$ws.Range("L86").Value = 6499.5
$ws.Range("M86").Value = -2876
$ws.Range("N86").Value = -8745.5
$ws.Range("H89").Value = 5249.25
$ws.Range("I89").Value = 3999
$ws.Range("J89").Value = 6499.5
$ws.Range("K89").Value = 19995
$ws.Range("L89").Value = 32497.5
$ws.Range("M89").Value = -14379
$ws.Range("N89").Value = -43729.5
$ws.Range("H99").Value = 61407.47
$ws.Range("I99").Value = 73316.21000000001
$ws.Range("K99").Value = 73316.21000000001
$ws.Range("M99").Value = -71818.21000000001
$ws.Range("H134").Value = 1157.125
$ws.Range("I134").Value = 729.0294
$ws.Range("K134").Value = 2187.0882
$ws.Range("M134").Value = 347.9117999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2361765.2
$ws.Range("I31").Value = 1440.963
$ws.Range("J31").Value = 4812871
$ws.Range("K31").Value = 1440.963
$ws.Range("L31").Value = 4812871
$ws.Range("M31").Value = -1145.963
$ws.Range("N31").Value = -4813461
$ws.Range("H34").Value = 2361765.2
$ws.Range("I34").Value = 1440.963
$ws.Range("J34").Value = 4812871
$ws.Range("K34").Value = 1440.963
$ws.Range("L34").Value = 4812871
$ws.Range("M34").Value = -1238.963
$ws.Range("N34").Value = -4813275
$ws.Range("H99").Value = 3191.6667
$ws.Range("I99").Value = 1974.1666
$ws.Range("K99").Value = 1974.1666
$ws.Range("M99").Value = -476.1666
$ws.Range("H126").Value = 3191.6667
$ws.Range("I126").Value = 1974.1666
$ws.Range("K126").Value = 5922.4998
$ws.Range("M126").Value = -3452.4998
$ws.Range("H134").Value = 3475.5
$ws.Range("I134").Value = 3572.1191
$ws.Range("K134").Value = 10716.3573
$ws.Range("M134").Value = -8181.3573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3662.3333
$ws.Range("I80").Value = 2993.5
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 8980.5
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -8044.5
$ws.Range("N80").Value = -16872
$ws.Range("H83").Value = 3662.3333
$ws.Range("I83").Value = 2993.5
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 26941.5
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -22261.5
$ws.Range("N83").Value = -54360
$ws.Range("H139").Value = 5208.3447
$ws.Range("I139").Value = 2255.25
$ws.Range("K139").Value = 6765.75
$ws.Range("M139").Value = -1625.75
$ws.Range("H140").Value = 2271.8635
$ws.Range("I140").Value = 2271.8635
$ws.Range("K140").Value = 6815.5905
$ws.Range("M140").Value = -1635.5905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5703.4287
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 5820.6665
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 17461.9995
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -22361.9995
$ws.Range("H132").Value = 2928.75
$ws.Range("I132").Value = 2643.5
$ws.Range("K132").Value = 7930.5
$ws.Range("M132").Value = -5400.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 8498.333000000001
$ws.Range("I26").Value = 8498.333000000001
$ws.Range("K26").Value = 8498.333000000001
$ws.Range("M26").Value = -8203.333000000001
$ws.Range("H100").Value = 6072.364
$ws.Range("I100").Value = 5711.1113
$ws.Range("K100").Value = 5711.1113
$ws.Range("M100").Value = -5170.1113
$ws.Range("H136").Value = 5477.222
$ws.Range("I136").Value = 3459.2
$ws.Range("K136").Value = 10377.6
$ws.Range("M136").Value = -7827.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2994.8
$ws.Range("I62").Value = 2487.5
$ws.Range("J62").Value = 3333
$ws.Range("K62").Value = 2487.5
$ws.Range("L62").Value = 3333
$ws.Range("M62").Value = -1863.5
$ws.Range("N62").Value = -4581
$ws.Range("H65").Value = 2994.8
$ws.Range("I65").Value = 2487.5
$ws.Range("J65").Value = 3333
$ws.Range("K65").Value = 12437.5
$ws.Range("L65").Value = 16665
$ws.Range("M65").Value = -9317.5
$ws.Range("N65").Value = -22905

